$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13 ("2021年") continues the yearly series that ends at row 12 ("2020年").
$row = 13
$prevRow = 12

# Column A ("year" label) - carries the same bold/centered/bordered style as the
# other year cells (A2:A12), so copy that formatting across before setting the value.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = "2021年"

# Numeric columns B..S.
$ws.Cells.Item($row, 2).Value  = 3427.1
$ws.Cells.Item($row, 3).Value  = 9406.9
$ws.Cells.Item($row, 4).Value  = 82547.8

# E and M have no reported figure for this year (same as rows 5-12) - an empty
# string, not a blank cell.
$ws.Cells.Item($row, 5).Formula  = '=""'
$ws.Cells.Item($row, 13).Formula = '=""'

$ws.Cells.Item($row, 6).Value  = 447330.3
$ws.Cells.Item($row, 7).Value  = 640711.3
$ws.Cells.Item($row, 8).Value  = 69804.3
$ws.Cells.Item($row, 9).Value  = 36611.9
$ws.Cells.Item($row, 10).Value = 60761.4
$ws.Cells.Item($row, 11).Value = 129184.2
$ws.Cells.Item($row, 12).Value = 593725.6

$ws.Cells.Item($row, 14).Value = 11727.1
$ws.Cells.Item($row, 15).Value = 42159.5
$ws.Cells.Item($row, 16).Value = 81058.89999999999
$ws.Cells.Item($row, 17).Value = 722909.8
$ws.Cells.Item($row, 18).Value = 842.2
$ws.Cells.Item($row, 19).Value = 30077.4
